$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.369.39"
$ws.Range("E2").Value = "  -0.61%  "

$ws.Range("D3").Value = "2.390.91"
$ws.Range("E3").Value = "  -3.86%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'549.73"
$ws.Range("E5").Value = "  -1.18%  "

$ws.Range("D6").Value = "'141.79"
$ws.Range("E6").Value = "  -4.50%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("E8").Value = "  -11.38%  "

$ws.Range("D9").Value = "2.389.67"
$ws.Range("E9").Value = "  -3.90%  "

$ws.Range("E10").Value = "  -3.02%  "

$ws.Range("E11").Value = "  +0.19%  "

$ws.Range("E12").Value = "  -3.58%  "

$ws.Range("E13").Value = "  -3.18%  "

$ws.Range("D14").Value = "'25.51"
$ws.Range("E14").Value = "  -3.98%  "

$ws.Range("D15").Value = "2.822.66"
$ws.Range("E15").Value = "  -3.87%  "

$ws.Range("D17").Value = "60.968.41"
$ws.Range("E17").Value = "  -1.08%  "

$ws.Range("D18").Value = "2.388.86"
$ws.Range("E18").Value = "  -3.88%  "

$ws.Range("E20").Value = "  -2.80%  "

$ws.Range("D21").Value = "'319.22"
$ws.Range("E21").Value = "  -1.34%  "

$ws.Range("D22").Value = "'6.71"
$ws.Range("E22").Value = "  -7.00%  "

$ws.Range("E23").Value = "  -0.02%  "

$ws.Range("E24").Value = "  -0.14%  "

$ws.Range("D25").Value = "'63.88"

$ws.Range("E26").Value = "  +3.47%  "

$ws.Range("D27").Value = "'0.998"
$ws.Range("E27").Value = "  -0.09%  "

$ws.Range("D28").Value = "2.508.34"
$ws.Range("E28").Value = "  -3.64%  "

$ws.Range("D33").Value = "'0.147"
$ws.Range("E33").Value = "  -4.33%  "

$ws.Range("D34").Value = "'1.84"
$ws.Range("E34").Value = "  -4.35%  "

$ws.Range("E35").Value = "  -1.34%  "

$ws.Range("E36").Value = "  -0.01%  "

$ws.Range("D37").Value = "'5.52"
$ws.Range("E37").Value = "  -7.77%  "

$ws.Range("D38").Value = "'4.68"
$ws.Range("E38").Value = "  -5.57%  "

$ws.Range("E39").Value = "  -2.55%  "

$ws.Range("D40").Value = "'1.85"
$ws.Range("E40").Value = "  +3.52%  "

$ws.Range("D41").Value = "'18.13"
$ws.Range("E41").Value = "  -3.03%  "

$ws.Range("D42").Value = "'140.03"
$ws.Range("E42").Value = "  -4.81%  "

$ws.Range("E43").Value = "  +0.05%  "

$ws.Range("D44").Value = "'40.42"
$ws.Range("E44").Value = "  -0.48%  "

$ws.Range("E45").Value = "  -11.39%  "

$ws.Range("E46").Value = "  -1.46%  "

$ws.Range("D47").Value = "'140.92"
$ws.Range("E47").Value = "  -5.42%  "

$ws.Range("D48").Value = "'20.20"
$ws.Range("E48").Value = "  -9.04%  "

$ws.Range("D49").Value = "'0.0521"
$ws.Range("E49").Value = "  -4.57%  "

$ws.Range("E50").Value = "  -3.86%  "

$ws.Range("E51").Value = "  -4.30%  "

$ws.Range("B29").Value = "Bittensor"
$ws.Range("C29").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D29").Value = "'530.98"
$ws.Range("E29").Value = "  -4.97%  "

$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0₃0931"
$ws.Range("E30").Value = "  -9.03%  "

$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "'8.13"
$ws.Range("E31").Value = "  -3.61%  "

$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "'1.43"
$ws.Range("E32").Value = "  -6.05%  "
